# Refresh crypto price/volume snapshot (GitHub Actions data pull).
# Column D ("Price") holds numeric-looking text (e.g. "59.538.97", "0.0000138")
# that must stay plain text, exactly as authored in the sheet (t="inlineStr"/"s",
# General format, no thousands grouping). A leading apostrophe forces Excel to
# store the value as text instead of silently parsing/rounding it as a number
# (and Excel strips the apostrophe itself, so the stored value is unaffected).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''59.538.97'
$ws.Range("E2").Value = '  -2.19%  '
# Row 3
$ws.Range("D3").Value = '''2.592.22'
$ws.Range("E3").Value = '  -2.05%  '
# Row 4
$ws.Range("E4").Value = '  +0.13%  '
# Row 5
$ws.Range("D5").Value = '''561.63'
$ws.Range("E5").Value = '  -1.07%  '
# Row 6
$ws.Range("D6").Value = '''143.64'
$ws.Range("E6").Value = '  -2.36%  '
# Row 7
$ws.Range("E7").Value = '  +0.20%  '
# Row 8
$ws.Range("E8").Value = '  -2.31%  '
# Row 9
$ws.Range("D9").Value = '''2.602.82'
$ws.Range("E9").Value = '  -2.50%  '
# Row 11
$ws.Range("D11").Value = '''0.105'
$ws.Range("E11").Value = '  -0.34%  '
# Row 12
$ws.Range("D12").Value = '''0.160'
$ws.Range("E12").Value = '  +10.78%  '
# Row 13
$ws.Range("E13").Value = '  +4.28%  '
# Row 14
$ws.Range("D14").Value = '''3.050.75'
$ws.Range("E14").Value = '  -2.25%  '
# Row 15
$ws.Range("D15").Value = '''23.41'
$ws.Range("E15").Value = '  +6.99%  '
# Row 16
$ws.Range("D16").Value = '''59.463.91'
$ws.Range("E16").Value = '  -1.74%  '
# Row 17
$ws.Range("D17").Value = '''0.0000138'
$ws.Range("E17").Value = '  +0.31%  '
# Row 18
$ws.Range("D18").Value = '''2.594.14'
$ws.Range("E18").Value = '  -3.04%  '
# Row 19
$ws.Range("D19").Value = '''4.60'
$ws.Range("E19").Value = '  +0.79%  '
# Row 20
$ws.Range("D20").Value = '''338.88'
$ws.Range("E20").Value = '  -1.41%  '
# Row 21
$ws.Range("D21").Value = '''10.44'
$ws.Range("E21").Value = '  -0.27%  '
# Row 22
$ws.Range("D22").Value = '''6.58'
$ws.Range("E22").Value = '  +3.23%  '
# Row 23
$ws.Range("E23").Value = '  +0.37%  '
# Row 24
$ws.Range("D24").Value = '''63.88'
$ws.Range("E24").Value = '  -4.09%  '
# Row 25
$ws.Range("D25").Value = '''0.471'
$ws.Range("E25").Value = '  +7.00%  '
# Row 26
$ws.Range("E26").Value = '  +0.36%  '
# Row 27
$ws.Range("D27").Value = '''0.162'
$ws.Range("E27").Value = '  -2.01%  '
# Row 28
$ws.Range("D28").Value = '''7.48'
$ws.Range("E28").Value = '  +1.12%  '
# Row 29
$ws.Range("D29").Value = '0.0₃0783'
$ws.Range("E29").Value = '  +0.06%  '
# Row 30
$ws.Range("E30").Value = '  +0.06%  '
# Row 31
$ws.Range("D31").Value = '''6.20'
$ws.Range("E31").Value = '  -2.86%  '
# Row 32
$ws.Range("D32").Value = '''1.68'
$ws.Range("E32").Value = '  -2.14%  '
# Row 33
$ws.Range("D33").Value = '''158.28'
$ws.Range("E33").Value = '  +2.56%  '
# Row 34
$ws.Range("D34").Value = '''19.13'
$ws.Range("E34").Value = '  -0.49%  '
# Row 35
$ws.Range("D35").Value = '''4.07'
$ws.Range("E35").Value = '  -0.45%  '
# Row 36
$ws.Range("E36").Value = '  +1.09%  '
# Row 37
$ws.Range("D37").Value = '''0.896'
$ws.Range("E37").Value = '  -0.96%  '
# Row 38
$ws.Range("D38").Value = '''0.876'
$ws.Range("E38").Value = '  -3.31%  '
# Row 39
$ws.Range("D39").Value = '''37.46'
# Row 40
$ws.Range("E40").Value = '  -1.46%  '
# Row 41
$ws.Range("D41").Value = '''3.69'
$ws.Range("E41").Value = '  +0.80%  '
# Row 42
$ws.Range("D42").Value = '''294.80'
$ws.Range("E42").Value = '  -2.59%  '
# Row 43
$ws.Range("D43").Value = '''139.64'
$ws.Range("E43").Value = '  +8.99%  '
# Row 44
$ws.Range("D44").Value = '''0.999'
$ws.Range("E44").Value = '  +0.18%  '
# Row 45
$ws.Range("D45").Value = '''0.0977'
$ws.Range("E45").Value = '  -0.57%  '
# Row 46
$ws.Range("E46").Value = '  -1.56%  '
# Row 47
$ws.Range("B47").Value = 'WhiteBITCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D47").Value = '''10.64'
$ws.Range("E47").Value = '  -0.27%  '
# Row 48
$ws.Range("B48").Value = 'Hedera'
$ws.Range("C48").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D48").Value = '''0.0531'
$ws.Range("E48").Value = '  -3.21%  '
# Row 49
$ws.Range("D49").Value = '''0.0236'
$ws.Range("E49").Value = '  +0.26%  '
# Row 50
$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D50").Value = '''4.71'
$ws.Range("E50").Value = '  +1.31%  '
# Row 51
$ws.Range("B51").Value = 'InjectiveProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D51").Value = '''18.87'
$ws.Range("E51").Value = '  -0.53%  '
